$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.262.90'
$ws.Range("E2").Value = '  +1.77%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.390.23'
$ws.Range("E3").Value = '  +7.46%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.22'
$ws.Range("E5").Value = '  +10.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.03'
$ws.Range("E6").Value = '  -5.65%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.650'
$ws.Range("E7").Value = '  +4.26%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.647'
$ws.Range("E9").Value = '  +7.76%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.09'
$ws.Range("E10").Value = '  -3.70%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0940'
$ws.Range("E11").Value = '  +2.87%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.58'
$ws.Range("E12").Value = '  -0.71%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.34'
$ws.Range("E13").Value = '  +16.03%  '

$ws.Range("E14").Value = '  -2.10%  '

$ws.Range("E15").Value = '  +2.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.752.93'
$ws.Range("E16").Value = '  +7.51%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.394.94'
$ws.Range("E17").Value = '  +7.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.269.22'
$ws.Range("E18").Value = '  +1.95%  '

$ws.Range("E19").Value = '  +3.70%  '

$ws.Range("E20").Value = '  +3.32%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '76.44'
$ws.Range("E21").Value = '  +3.84%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '273.72'
$ws.Range("E22").Value = '  +16.58%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.42'
$ws.Range("E23").Value = '  +2.55%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.41'
$ws.Range("E24").Value = '  +0.50%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.73'
$ws.Range("E25").Value = '  +9.43%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.83'
$ws.Range("E26").Value = '  +3.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.08%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.96'
$ws.Range("E28").Value = '  +7.99%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '177.22'
$ws.Range("E29").Value = '  +1.58%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.00'
$ws.Range("E30").Value = '  +1.96%  '

$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.20'
$ws.Range("E31").Value = '  -0.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.22'
$ws.Range("E32").Value = '  +2.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0930'
$ws.Range("E33").Value = '  +5.50%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.89'
$ws.Range("E34").Value = '  +3.74%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.133'
$ws.Range("E35").Value = '  +6.13%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.87'
$ws.Range("E36").Value = '  -2.89%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.15'
$ws.Range("E37").Value = '  -0.82%  '

$ws.Range("E38").Value = '  -1.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.108'
$ws.Range("E39").Value = '  +3.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.84'
$ws.Range("E40").Value = '  +18.92%  '

$ws.Range("E41").Value = '  +21.77%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '125.51'
$ws.Range("E42").Value = '  +24.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.233'
$ws.Range("E43").Value = '  +1.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '69.19'
$ws.Range("E44").Value = '  -3.33%  '

$ws.Range("E45").Value = '  +0.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.56'
$ws.Range("E46").Value = '  +1.83%  '

$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '92.19'
$ws.Range("E47").Value = '  +68.27%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.55'
$ws.Range("E48").Value = '  +13.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.64'
$ws.Range("E49").Value = '  +5.18%  '

$ws.Range("E50").Value = '  +2.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.598.59'
$ws.Range("E51").Value = '  +11.77%  '
